{"js": "// The body of: async (context) => { ... }\n\nconst body = context.document.body;\n\n// ------------------------------------------------------------------\n// 1) \"...power consumption -- will depend on the on the particular\n//    application.\" becomes \"...power consumption etc -- will depend\n//    on the on the particular application being considered.\"\n//\n//    This is done as two separate, run-boundary-preserving edits\n//    (one per original run) instead of one big replace spanning both\n//    runs, so unrelated neighboring runs are not disturbed/merged.\n// ------------------------------------------------------------------\nconst consumptionResults = body.search(\", power consumption \", {\n  matchCase: true,\n});\nconsumptionResults.load(\"text\");\nawait context.sync();\n\nif (consumptionResults.items.length > 0) {\n  consumptionResults.items[0].insertText(\n    \", power consumption etc \",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\nconst dependResults = body.search(\n  \"-- will depend on the on the particular application.\",\n  { matchCase: true }\n);\ndependResults.load(\"text\");\nawait context.sync();\n\nif (dependResults.items.length > 0) {\n  dependResults.items[0].insertText(\n    \"-- will depend on the on the particular application being considered.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// ------------------------------------------------------------------\n// 2) Move the hidden \"_GoBack\" bookmark from right before \"We choose a\n//    relatively simple impulse response\" to right before the final\n//    period that ends \"...outside the scope of this article.\"\n// ------------------------------------------------------------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst articleResults = body.search(\"outside the scope of this article\", {\n  matchCase: true,\n});\narticleResults.load(\"text\");\nawait context.sync();\n\nif (articleResults.items.length > 0) {\n  const endOfArticle = articleResults.items[0].getRange(Word.RangeLocation.end);\n  endOfArticle.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// ------------------------------------------------------------------\n// 3) Merge the trailing '\" will depend...future studies.' run with the\n//    closing quote run ('\u201d ') into a single run. There is no text\n//    change here, so re-writing the combined range in place causes the\n//    two adjacent runs (identical formatting) to collapse into one.\n// ------------------------------------------------------------------\nconst closingResults = body.search(\n  \" while detailed studies of more realistic circuit implementations and correlated effects are left for future studies.\\u201D \",\n  { matchCase: true }\n);\nclosingResults.load(\"text\");\nawait context.sync();\n\nif (closingResults.items.length > 0) {\n  closingResults.items[0].insertText(\n    \" while detailed studies of more realistic circuit implementations and correlated effects are left for future studies.\\u201D \",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is pre-seeded as $d below.\n\n$d = $word.ActiveDocument\n\n# ------------------------------------------------------------------\n# 1) \"...power consumption -- will depend on the on the particular\n#    application.\" becomes \"...power consumption etc -- will depend\n#    on the on the particular application being considered.\"\n# ------------------------------------------------------------------\n$oldSentence = \", power consumption -- will depend on the on the particular application.\"\n$newSentence = \", power consumption etc -- will depend on the on the particular application being considered.\"\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Execute($oldSentence, $false, $false, $false, $false, $false, $true, 1, $false, $newSentence, 2) | Out-Null\n\n# ------------------------------------------------------------------\n# 2) Move the hidden \"_GoBack\" bookmark from right before \"We choose a\n#    relatively simple impulse response\" to right before the final\n#    period that ends \"...outside the scope of this article.\"\n# ------------------------------------------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$found2 = $rng2.Find.Execute(\"outside the scope of this article\")\nif ($found2) {\n    $endOfArticle = $rng2.Duplicate\n    $endOfArticle.Collapse(0)   # wdCollapseEnd\n    $d.Bookmarks.Add(\"_GoBack\", $endOfArticle) | Out-Null\n}\n\n# ------------------------------------------------------------------\n# 3) Merge the trailing \"...future studies.\" run with the closing\n#    quote run (\"\\u201d \") into a single run. There is no text change\n#    here; re-running Find & Replace with identical text forces Word\n#    to collapse the two adjacent, identically-formatted runs.\n# ------------------------------------------------------------------\n$closingText = \" while detailed studies of more realistic circuit implementations and correlated effects are left for future studies.\" + [char]8221 + \" \"\n\n$rng3 = $d.Content\n$rng3.Find.ClearFormatting()\n$rng3.Find.Execute($closingText, $false, $false, $false, $false, $false, $true, 1, $false, $closingText, 2) | Out-Null\n"}
